$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Reference")

$ws.Range("M1:V15").EntireColumn.Delete()
